# "add completed event when we open notepad"
# Adds a new `notepad` key (two lines: notepad_1 / notepad_2) to the
# keys sheet plus its en-US and ru-RU translations, bumps the "count of
# keys" counter on the main sheet, and leaves the workbook focused on the
# ru-RU sheet (matching the author's final selection/view state).

$wb = $excel.ActiveWorkbook

# --- main sheet: bump "count of keys" from 8 to 10 -------------------
$main = $wb.Worksheets.Item("main")
$main.Range("B2").Value = 10

# --- keys sheet: register the two new translation keys ----------------
$keys = $wb.Worksheets.Item("keys")
$keys.Range("A9").Value = "notepad_1"
$keys.Range("B9").Value = 4
$keys.Range("A10").Value = "notepad_2"
$keys.Range("B10").Value = 4
$keys.Range("B10").Select()

# --- ru-RU sheet: Russian translations ---------------------------------
$ruRU = $wb.Worksheets.Item("ru-RU")
$ruRU.Range("A9").Value = "notepad_1"
$ruRU.Range("B9").Value = "Эй, ты открыл блокнот?"
$ruRU.Range("A10").Value = "notepad_2"
$ruRU.Range("B10").Value = "А зачем?"

# --- en-US sheet: English translations ----------------------------------
$enUS = $wb.Worksheets.Item("en-US")
$enUS.Range("A9").Value = "notepad_1"
$enUS.Range("B9").Value = "Do you open notepad?"
$enUS.Range("A10").Value = "notepad_2"
$enUS.Range("B10").Value = "But why?"
$enUS.Range("A9:B10").HorizontalAlignment = -4108
$enUS.Range("A9:B10").VerticalAlignment = -4108
$enUS.Range("B10").Select()

# ru-RU ends up the active/selected sheet, with C19 selected
$ruRU.Activate()
$ruRU.Range("C19").Select()
